# Auto-generated edit script applying numeric corrections per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 2710
$ws.Range("J52").Value = 2710
$ws.Range("L52").Value = 8130
$ws.Range("N52").Value = -8450
$ws.Range("H58").Value = 1781.8572
$ws.Range("I58").Value = 158.33333
$ws.Range("J58").Value = 2999.5
$ws.Range("K58").Value = 474.99999
$ws.Range("L58").Value = 8998.5
$ws.Range("M58").Value = -324.99999
$ws.Range("N58").Value = -9298.5
$ws.Range("H64").Value = 458933.62
$ws.Range("I64").Value = 717817.2
$ws.Range("J64").Value = 5887.375
$ws.Range("K64").Value = 717817.2
$ws.Range("L64").Value = 5887.375
$ws.Range("M64").Value = -717569.2
$ws.Range("N64").Value = -6383.375
$ws.Range("H67").Value = 458933.62
$ws.Range("I67").Value = 717817.2
$ws.Range("J67").Value = 5887.375
$ws.Range("K67").Value = 717817.2
$ws.Range("L67").Value = 5887.375
$ws.Range("M67").Value = -716959.2
$ws.Range("N67").Value = -7603.375
$ws.Range("H70").Value = 2140.3635
$ws.Range("I70").Value = 997.8
$ws.Range("J70").Value = 3092.5
$ws.Range("K70").Value = 2993.4
$ws.Range("L70").Value = 9277.5
$ws.Range("M70").Value = -2723.4
$ws.Range("N70").Value = -9817.5
$ws.Range("H73").Value = 2140.3635
$ws.Range("I73").Value = 997.8
$ws.Range("J73").Value = 3092.5
$ws.Range("K73").Value = 2993.4
$ws.Range("L73").Value = 9277.5
$ws.Range("M73").Value = -2057.4
$ws.Range("N73").Value = -11149.5
$ws.Range("H112").Value = 6199538.5
$ws.Range("J112").Value = 6199538.5
$ws.Range("L112").Value = 18598615.5
$ws.Range("N112").Value = -18600831.5
$ws.Range("H138").Value = 5935284.5
$ws.Range("I138").Value = 1918520.1
$ws.Range("J138").Value = 7578506
$ws.Range("K138").Value = 5755560.300000001
$ws.Range("L138").Value = 22735518
$ws.Range("M138").Value = -5750420.300000001
$ws.Range("N138").Value = -22745798

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 322
$ws.Range("I4").Value = 322
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 322
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = -206
$ws.Range("H5").Value = 236.14285
$ws.Range("I5").Value = 230.2
$ws.Range("J5").Value = 251
$ws.Range("K5").Value = 230.2
$ws.Range("L5").Value = 251
$ws.Range("M5").Value = -118.2
$ws.Range("N5").Value = -475
$ws.Range("H74").Value = 3932.8
$ws.Range("I74").Value = 1235.2424
$ws.Range("J74").Value = 11351.083
$ws.Range("K74").Value = 1235.2424
$ws.Range("L74").Value = 11351.083
$ws.Range("M74").Value = -361.2424000000001
$ws.Range("N74").Value = -13099.083
$ws.Range("H77").Value = 3932.8
$ws.Range("I77").Value = 1235.2424
$ws.Range("J77").Value = 11351.083
$ws.Range("K77").Value = 6176.212
$ws.Range("L77").Value = 56755.415
$ws.Range("M77").Value = -1808.212
$ws.Range("N77").Value = -65491.415
$ws.Range("M4").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 236.14285
$ws.Range("I4").Value = 230.2
$ws.Range("J4").Value = 251
$ws.Range("K4").Value = 230.2
$ws.Range("L4").Value = 251
$ws.Range("M4").Value = -115.2
$ws.Range("N4").Value = -481
$ws.Range("H32").Value = 50000
$ws.Range("I32").Value = 50000
$ws.Range("K32").Value = 50000
$ws.Range("M32").Value = -49616
$ws.Range("H134").Value = 3192.484
$ws.Range("I134").Value = 2197.8333
$ws.Range("J134").Value = 6602.7144
$ws.Range("K134").Value = 6593.499899999999
$ws.Range("L134").Value = 19808.1432
$ws.Range("M134").Value = -4058.499899999999
$ws.Range("N134").Value = -24878.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 123
$ws.Range("I7").Value = 84.5
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 84.5
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = 28.5
$ws.Range("N7").Value = -426
$ws.Range("H59").Value = 27162.5
$ws.Range("J59").Value = 27162.5
$ws.Range("L59").Value = 27162.5
$ws.Range("N59").Value = -29452.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 19231678
$ws.Range("J113").Value = 21740080
$ws.Range("L113").Value = 65220240
$ws.Range("N113").Value = -65224580
$ws.Range("H132").Value = 898.8570999999999
$ws.Range("I132").Value = 674.25
$ws.Range("J132").Value = 1198.3334
$ws.Range("K132").Value = 6068.25
$ws.Range("L132").Value = 10785.0006
$ws.Range("M132").Value = -3538.25
$ws.Range("N132").Value = -15845.0006
$ws.Range("H136").Value = 4827.864
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 5160.65
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 15481.95
$ws.Range("M136").Value = 600
$ws.Range("N136").Value = -25681.95
$ws.Range("H139").Value = 1923.1538
$ws.Range("I139").Value = 1655.174
$ws.Range("J139").Value = 3977.6667
$ws.Range("K139").Value = 4965.522
$ws.Range("L139").Value = 11933.0001
$ws.Range("M139").Value = 174.4780000000001
$ws.Range("N139").Value = -22213.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10574.8
$ws.Range("J123").Value = 10574.8
$ws.Range("L123").Value = 10574.8
$ws.Range("N123").Value = -15474.8
$ws.Range("H139").Value = 32549.5
$ws.Range("J139").Value = 32549.5
$ws.Range("L139").Value = 32549.5
$ws.Range("N139").Value = -42829.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3639.7925
$ws.Range("I132").Value = 2900.8845
$ws.Range("J132").Value = 4351.3335
$ws.Range("K132").Value = 8702.6535
$ws.Range("L132").Value = 13054.0005
$ws.Range("M132").Value = -6172.6535
$ws.Range("N132").Value = -18114.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 40127
$ws.Range("I122").Value = 54105.42
$ws.Range("J122").Value = 2185.5715
$ws.Range("K122").Value = 162316.26
$ws.Range("L122").Value = 6556.7145
$ws.Range("M122").Value = -159866.26
$ws.Range("N122").Value = -11456.7145
$ws.Range("H126").Value = 37273.43
$ws.Range("I126").Value = 49095.285
$ws.Range("J126").Value = 1807.8572
$ws.Range("K126").Value = 147285.855
$ws.Range("L126").Value = 5423.571599999999
$ws.Range("M126").Value = -144815.855
$ws.Range("N126").Value = -10363.5716
